$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "51.060.03"
$cell.Style = $oldStyle
$ws.Range("E2").Value = "  +0.05%  "

$cell = $ws.Range("D3")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.957.80"
$cell.Style = $oldStyle
$ws.Range("E3").Value = "  +0.71%  "

$cell = $ws.Range("D4")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = $oldStyle
$ws.Range("E4").Value = "  +0.02%  "

$cell = $ws.Range("D5")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "380.17"
$cell.Style = $oldStyle
$ws.Range("E5").Value = "  +1.78%  "

$cell = $ws.Range("D6")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "101.82"
$cell.Style = $oldStyle
$ws.Range("E6").Value = "  -0.41%  "

$cell = $ws.Range("D7")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.544"
$cell.Style = $oldStyle
$ws.Range("E7").Value = "  +1.67%  "

$ws.Range("E8").Value = "  -0.04%  "

$cell = $ws.Range("D9")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.586"
$cell.Style = $oldStyle
$ws.Range("E9").Value = "  +0.78%  "

$cell = $ws.Range("D10")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "36.66"
$cell.Style = $oldStyle
$ws.Range("E10").Value = "  +0.24%  "

$ws.Range("E11").Value = "  -1.06%  "

$ws.Range("E12").Value = "  +2.30%  "

$cell = $ws.Range("D13")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.422.36"
$cell.Style = $oldStyle
$ws.Range("E13").Value = "  +0.69%  "

$cell = $ws.Range("D14")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "18.33"
$cell.Style = $oldStyle
$ws.Range("E14").Value = "  +2.51%  "

$ws.Range("E15").Value = "  +5.66%  "

$cell = $ws.Range("D16")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "12.26"
$cell.Style = $oldStyle
$ws.Range("E16").Value = "  +72.29%  "

$cell = $ws.Range("D17")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.969.67"
$cell.Style = $oldStyle
$ws.Range("E17").Value = "  +1.04%  "

$ws.Range("E18").Value = "  +2.98%  "

$cell = $ws.Range("D19")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "51.130.05"
$cell.Style = $oldStyle
$ws.Range("E19").Value = "  +0.24%  "

$cell = $ws.Range("D20")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.09"
$cell.Style = $oldStyle
$ws.Range("E20").Value = "  -1.53%  "

$cell = $ws.Range("D21")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "12.38"
$cell.Style = $oldStyle
$ws.Range("E21").Value = "  -1.30%  "

$cell = $ws.Range("D22")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0961"
$cell.Style = $oldStyle
$ws.Range("E22").Value = "  +0.98%  "

$ws.Range("E23").Value = "  +15.57%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$cell = $ws.Range("D24")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "69.79"
$cell.Style = $oldStyle
$ws.Range("E24").Value = "  +2.39%  "

$ws.Range("B25").Value = "BitcoinCash"
$ws.Range("C25").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$cell = $ws.Range("D25")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "268.27"
$cell.Style = $oldStyle
$ws.Range("E25").Value = "  +1.81%  "

$cell = $ws.Range("D26")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "7.92"
$cell.Style = $oldStyle
$ws.Range("E26").Value = "  -1.88%  "

$cell = $ws.Range("D27")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = $oldStyle
$ws.Range("E27").Value = "  -0.03%  "

$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$cell = $ws.Range("D28")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "7.12"
$cell.Style = $oldStyle
$ws.Range("E28").Value = "  -9.76%  "

$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$cell = $ws.Range("D29")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.167"
$cell.Style = $oldStyle
$ws.Range("E29").Value = "  -0.79%  "

$cell = $ws.Range("D30")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "25.84"
$cell.Style = $oldStyle
$ws.Range("E30").Value = "  +0.80%  "

$cell = $ws.Range("D31")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.109"
$cell.Style = $oldStyle
$ws.Range("E31").Value = "  -3.40%  "

$cell = $ws.Range("D32")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "10.42"
$cell.Style = $oldStyle
$ws.Range("E32").Value = "  +5.96%  "

$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$cell = $ws.Range("D33")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "34.29"
$cell.Style = $oldStyle
$ws.Range("E33").Value = "  +0.86%  "

$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$cell = $ws.Range("D34")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "51.04"
$cell.Style = $oldStyle
$ws.Range("E34").Value = "  +0.07%  "

$ws.Range("B35").Value = "Toncoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$cell = $ws.Range("D35")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.07"
$cell.Style = $oldStyle
$ws.Range("E35").Value = "  +2.28%  "

$cell = $ws.Range("D36")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0435"
$cell.Style = $oldStyle
$ws.Range("E36").Value = "  -4.49%  "

$ws.Range("E37").Value = "  +0.05%  "

$ws.Range("E38").Value = "  +10.55%  "

$ws.Range("E39").Value = "  +1.94%  "

$cell = $ws.Range("D40")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "16.61"
$cell.Style = $oldStyle
$ws.Range("E40").Value = "  +1.51%  "

$ws.Range("E41").Value = "  +3.39%  "

$cell = $ws.Range("D42")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.49"
$cell.Style = $oldStyle
$ws.Range("E42").Value = "  -2.24%  "

$cell = $ws.Range("D43")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "124.03"
$cell.Style = $oldStyle
$ws.Range("E43").Value = "  +2.68%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$cell = $ws.Range("D44")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "21.60"
$cell.Style = $oldStyle
$ws.Range("E44").Value = "  +3.12%  "

$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$cell = $ws.Range("D45")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.54"
$cell.Style = $oldStyle
$ws.Range("E45").Value = "  +10.20%  "

$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$cell = $ws.Range("D46")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.40"
$cell.Style = $oldStyle
$ws.Range("E46").Value = "  +4.59%  "

$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$cell = $ws.Range("D47")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.02"
$cell.Style = $oldStyle
$ws.Range("E47").Value = "  -1.25%  "

$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$cell = $ws.Range("D48")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.057.24"
$cell.Style = $oldStyle
$ws.Range("E48").Value = "  +3.41%  "

$ws.Range("B49").Value = "TheGraph"
$ws.Range("C49").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$cell = $ws.Range("D49")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.268"
$cell.Style = $oldStyle
$ws.Range("E49").Value = "  -2.13%  "

$cell = $ws.Range("D50")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0318"
$cell.Style = $oldStyle
$ws.Range("E50").Value = "  -7.77%  "

$cell = $ws.Range("D51")
$oldStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.42"
$cell.Style = $oldStyle
$ws.Range("E51").Value = "  +7.75%  "
